# Add a "Forearms" section to the exercise list, inserted right after the
# "Biceps" rows and before the "Abs/Core" rows (i.e. before the current row 65).
# This pushes the existing Abs/Core rows down by 3 (they keep their values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert three new blank rows starting at row 65; existing rows 65-70
# (the Abs/Core section) shift down to rows 68-73 automatically, keeping
# their original values intact.
$ws.Rows("65:67").Insert()

# Populate the newly inserted rows with the Forearms data.
$ws.Cells.Item(65, 1).Value = "Forearms"
$ws.Cells.Item(65, 2).Value = "Reverse Bar Bicep Curl"
$ws.Cells.Item(65, 3).Value = "Brachialis / Forearm"
$ws.Cells.Item(65, 4).Value = "EZ-Bar / Barbell"

$ws.Cells.Item(66, 1).Value = "Forearms"
$ws.Cells.Item(66, 2).Value = "Dumbbell Forearm Curl"
$ws.Cells.Item(66, 3).Value = "Forearms"
$ws.Cells.Item(66, 4).Value = "Dumbbells"

$ws.Cells.Item(67, 1).Value = "Forearms"
$ws.Cells.Item(67, 2).Value = "Cable Forearm Curl"
$ws.Cells.Item(67, 3).Value = "Forearms"
$ws.Cells.Item(67, 4).Value = "Cable Machine"
